$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "layerB"
$ws.Range("C1").Value = "minicolumn"
$ws.Range("D1").Value = "hypercolumn"
$ws.Range("E1").Value = "layerA"

$ws.Range("A2").Value = 0.2549039664327203
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4

$ws.Range("A3").Value = 0.35478181682618654
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 23

$ws.Range("A4").Value = 0.31535174109020064
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 23

$ws.Range("A5").Value = 0.3338125275154786
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 4

$ws.Range("A6").Value = 0.17964031035408046
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 23

$ws.Range("A7").Value = 0.34204593875625355
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 23

$ws.Range("A8").Value = 0.2013907069132586
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 4

$ws.Range("A9").Value = 0.21018811629290227
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 23

$ws.Range("A10").Value = 0.18787740513630202
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 23

$ws.Range("A11").Value = 0.3439294826399343
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 4

$ws.Range("A12").Value = 0.22280972317589728
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 23

$ws.Range("A13").Value = 0.31757591604666235
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 23

$ws.Range("A14").Value = 0.35567511231240484
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 4

$ws.Range("A15").Value = 0.33313995020973647
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 23

$ws.Range("A16").Value = 0.3663173663228573
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 23

$ws.Range("A17").Value = 0.2525305753031878
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 4

$ws.Range("A18").Value = 0.16100610350914504
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 23

$ws.Range("A19").Value = 0.40947719609112343
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 23

$ws.Range("A20").Value = 0.28246368024816876
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 4

$ws.Range("A21").Value = 0.26067272897491006
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 23

$ws.Range("A22").Value = 0.34510038998390813
$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 23

$ws.Range("A23").Value = 0.31788385954450166
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 4

$ws.Range("A24").Value = 0.2598068604771137
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 23

$ws.Range("A25").Value = 0.25798145794965643
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 23

$ws.Range("A26").Value = 0.14599034404361172
$ws.Range("B26").Value = 5
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 4

$ws.Range("A27").Value = 0.25958920068627805
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 23

$ws.Range("A28").Value = 0.22066999620609992
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 23

$ws.Range("A29").Value = 0.29037571506590487
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 4

$ws.Range("A30").Value = 0.24368783432210978
$ws.Range("B30").Value = 4
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 23

$ws.Range("A31").Value = 0.29702787597327907
$ws.Range("B31").Value = 5
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 23

$ws.Range("A32").Value = 0.27590253632600414
$ws.Range("B32").Value = 5
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 4

$ws.Range("A33").Value = 0.31479413452861094
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 23

$ws.Range("A34").Value = 0.29412898455986064
$ws.Range("B34").Value = 5
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 23

$ws.Range("A35").Value = 0.20356719395034573
$ws.Range("B35").Value = 5
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 4

$ws.Range("A36").Value = 0.24363653310247163
$ws.Range("B36").Value = 4
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 23

$ws.Range("A37").Value = 0.32173694856257873
$ws.Range("B37").Value = 5
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 23
